# "simplify steel description (remove RME)"
# The B2 cell (industrial steel mapping description) loses the "/RME"
# fragment from its 5th line, the cell switches to wrapped text with a
# taller row, and the sheet selection is left spanning B2:B11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuild the multi-line description without the "/RME" segment on the
# last line (was "18% S/LFM+CDN/RME/H:1", now "18% S/LFM+CDN/H:1").
$ws.Range("B2").Value = "16% CR/LFM+CDN/H:2`n29% CR+PC/LFM+CDN/H:1`n32% S/LFBR+CDN/H:1`n5% W/LWAL+CDN/H:1`n18% S/LFM+CDN/H:1"

# Wrap the text in the cell and grow the row to fit the 5 lines.
$ws.Range("B2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 256

# Leave the sheet's selection spanning the description column.
$ws.Range("B2:B11").Select()
